$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 223: add date label (text, avoid auto date-conversion via quote-prefix, then reset style) ---
$ws.Range("A223").Value = "'01-06-2021"
$ws.Range("A223").Style = "Normal"

# --- New row 223: numeric values ---
$ws.Range("B223").Value = -3504
$ws.Range("C223").Value = -2760
$ws.Range("D223").Value = 0
$ws.Range("E223").Value = -744
$ws.Range("F223").Value = -460
$ws.Range("G223").Value = -2604
$ws.Range("H223").Value = -2596
$ws.Range("I223").Value = 0
$ws.Range("J223").Value = -9
$ws.Range("K223").Value = 588
$ws.Range("L223").Value = -89
$ws.Range("M223").Value = 153
$ws.Range("N223").Value = -581
$ws.Range("O223").Value = 1105
$ws.Range("P223").Value = 232
$ws.Range("Q223").Value = 0
$ws.Range("R223").Value = 419
$ws.Range("S223").Value = -216
$ws.Range("T223").Value = 29
$ws.Range("U223").Value = 1325
$ws.Range("V223").Value = 121
$ws.Range("W223").Value = 1213
$ws.Range("X223").Value = -12
$ws.Range("Y223").Value = 3
$ws.Range("Z223").Value = 2167
$ws.Range("AA223").Value = 511
$ws.Range("AB223").Value = 1461
$ws.Range("AC223").Value = -46
$ws.Range("AD223").Value = 241

# --- Revisions to existing rows 218-222 ---
$ws.Range("F218").Value = -1878
$ws.Range("P218").Value = 1854
$ws.Range("R218").Value = 2294
$ws.Range("T218").Value = -152
$ws.Range("U218").Value = 1103
$ws.Range("V218").Value = 249
$ws.Range("W218").Value = 932
$ws.Range("Y218").Value = -73
$ws.Range("Z218").Value = 1586
$ws.Range("AA218").Value = 955
$ws.Range("AB218").Value = 72
$ws.Range("AD218").Value = 667
$ws.Range("F219").Value = -321
$ws.Range("P219").Value = -26
$ws.Range("R219").Value = -133
$ws.Range("T219").Value = 343
$ws.Range("U219").Value = 170
$ws.Range("V219").Value = -73
$ws.Range("W219").Value = 310
$ws.Range("Y219").Value = -62
$ws.Range("Z219").Value = 7445
$ws.Range("AA219").Value = 6158
$ws.Range("AB219").Value = 965
$ws.Range("AD219").Value = 398
$ws.Range("F220").Value = 2148
$ws.Range("P220").Value = -1344
$ws.Range("R220").Value = -918
$ws.Range("T220").Value = -196
$ws.Range("U220").Value = 900
$ws.Range("V220").Value = 490
$ws.Range("W220").Value = 494
$ws.Range("Y220").Value = -79
$ws.Range("Z220").Value = 1410
$ws.Range("AA220").Value = 619
$ws.Range("AB220").Value = 289
$ws.Range("AD220").Value = 567
$ws.Range("F221").Value = 1291
$ws.Range("P221").Value = -53
$ws.Range("R221").Value = 356
$ws.Range("U221").Value = 213
$ws.Range("V221").Value = 45
$ws.Range("W221").Value = 202
$ws.Range("Y221").Value = -28
$ws.Range("Z221").Value = 3826
$ws.Range("AA221").Value = 391
$ws.Range("AB221").Value = 739
$ws.Range("AD221").Value = 2752
$ws.Range("F222").Value = -1001
$ws.Range("K222").Value = 1082
$ws.Range("L222").Value = -7
$ws.Range("O222").Value = 1125
$ws.Range("P222").Value = -7294
$ws.Range("R222").Value = -7385
$ws.Range("U222").Value = 184
$ws.Range("V222").Value = 74
$ws.Range("W222").Value = 80
$ws.Range("Y222").Value = 43
$ws.Range("Z222").Value = 657
$ws.Range("AA222").Value = 300
$ws.Range("AB222").Value = 970
$ws.Range("AD222").Value = -540
